$d = $word.ActiveDocument

# --- Paragraph 2: "This week we had to do ... we had to make a presentation ..." ---
$d.Content.Find.Execute("This week we had to do", $true, $false, $false, $false, $false, $true, 1, $false, "This week groups had to do", 2)
$d.Content.Find.Execute("employee. I was in a group of four members and we came up with", $true, $false, $false, $false, $false, $true, 1, $false, "employee. The group was formed of four members and they came up with", 2)
$d.Content.Find.Execute("Employability Game, and we had to make", $true, $false, $false, $false, $false, $true, 1, $false, "Employability Game, and they had to make", 2)

# --- Paragraph 3: dice description ---
$d.Content.Find.Execute("the board using a dice, which shows", $true, $false, $false, $false, $false, $true, 1, $false, "the board using a six number dice, which shows", 2)

# --- "e.g" -> "e.g." ---
$d.Content.Find.Execute("from the dice, e.g dice shows 5 moves", $true, $false, $false, $false, $false, $true, 1, $false, "from the dice, e.g. dice shows 5 moves", 2)

Write-Output "done"
